$d = $word.ActiveDocument

# 1. Fix the trailing double space in the "Veel succes..." paragraph.
$d.Content.Find.Execute("Veel succes alvast vanmiddag.  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Veel succes alvast vanmiddag. ", 2)

# 2. Add a new paragraph after it with the contact info, matching the
#    formatting of the preceding paragraph (inherited automatically).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Altijd te bereiken op jamie.goegebeur@intecbrussel.be "
